$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 45996
$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 2
$ws.Range("G127").Value = 241.38
$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 67
$ws.Range("G227").Value = 9666.76
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("B229").Value = 63531
$ws.Range("E229").Value = 152.53
$ws.Range("F229").Value = 70
$ws.Range("G229").Value = 10043.6
$ws.Range("B230").Value = 57802
$ws.Range("E230").Value = 162.71
$ws.Range("F230").Value = -79
$ws.Range("G230").Value = -11334.92
$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12
$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 127
$ws.Range("G233").Value = 6050.28
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("B322").Value = 47097
$ws.Range("D322").Value = 112.28
$ws.Range("E322").Value = 134.16
$ws.Range("F322").Value = 15
$ws.Range("G322").Value = 1684.2
$ws.Range("B323").Value = 58047
$ws.Range("D323").Value = 105.54
$ws.Range("E323").Value = 126.1
$ws.Range("F323").Value = 41
$ws.Range("G323").Value = 4327.14
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29
$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9
$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945
$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65
$ws.Range("B382").Value = 64919
$ws.Range("E382").Value = 27.97
$ws.Range("F382").Value = 61
$ws.Range("G382").Value = 1604.3
$ws.Range("B383").Value = 45702
$ws.Range("E383").Value = 31.43
$ws.Range("F383").Value = -215
$ws.Range("G383").Value = -5654.5
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 5
$ws.Range("G442").Value = 1369.6
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34
